$wb = $excel.ActiveWorkbook

$wsFixedDiameter = $wb.Worksheets.Item("fixed_diameter")
$wsGeophysical = $wb.Worksheets.Item("geophysical")

# Update the Monte Carlo input values (rows 6-7) on both sheets.
$wsFixedDiameter.Range("E6").Value = 100
$wsFixedDiameter.Range("E7").Value = 50
$wsFixedDiameter.Range("F7").Value = 150

$wsGeophysical.Range("E6").Value = 100
$wsGeophysical.Range("E7").Value = 50
$wsGeophysical.Range("F7").Value = 150

# Update the selection shown on the "geophysical" sheet and leave it not
# as the active tab.
$wsGeophysical.Activate()
$wsGeophysical.Range("E6:F7").Select()

# Make "fixed_diameter" the active (selected) tab with its own selection,
# activating it last so it becomes the workbook's active sheet.
$wsFixedDiameter.Activate()
$wsFixedDiameter.Range("H13").Select()
